$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "10_291115_2231_7_xgboost_with_random_3in1_preprocess_valid1_valid2_"
$ws.Range("B11").Value = 0.62914999999999999
$ws.Range("C11").Value = "ensembled 7 xgboost, in random 3in1 data set with features preprocessed, with 2 valid sets"
